# Auto-generated: applies scheduled-runner price/profit refresh to the
# per-job Leve-flip sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Updated cells come from a fresh Universalis price pull; LeveProfit*
# columns (M/N) are recomputed from the new prices, including a couple
# of rows where a profit column is newly in/out of range.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H62" = 6691.769
    "I62" = 5332.1665
    "K62" = 5332.1665
    "M62" = -4708.1665
    "H65" = 6691.769
    "I65" = 5332.1665
    "K65" = 26660.8325
    "M65" = -23540.8325
    "H105" = 90922
    "J105" = 90922
    "L105" = 90922
    "N105" = -97910
    "H107" = 1592.4286
    "I107" = 899.5454999999999
    "K107" = 899.5454999999999
    "M107" = 1020.4545
    "H111" = 1000
    "I111" = 1000
    "K111" = 3000
    "M111" = 67
    "H132" = 959.7857
    "I132" = 976.8148
    "J132" = 500
    "K132" = 2930.4444
    "L132" = 1500
    "M132" = -400.4443999999999
    "N132" = -6560
    "H137" = 4443.8
    "I137" = 2929.8667
    "K137" = 8789.6001
    "M137" = -6239.6001
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H32" = 22729582
    "I32" = 25001640
    "K32" = 25001640
    "M32" = -25001353
    "H122" = 1388.2307
    "I122" = 1087.3334
    "K122" = 3262.0002
    "M122" = -812.0001999999999
    "H132" = 4957.9473
    "I132" = 2683.4644
    "J132" = 11326.5
    "K132" = 8050.3932
    "L132" = 33979.5
    "M132" = -5520.3932
    "N132" = -39039.5
    "H137" = 290000
    "J137" = 290000
    "L137" = 290000
    "N137" = -300200
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H62" = 0
    "I62" = 0
    "J62" = 0
    "K62" = 0
    "L62" = 0
    "H65" = 0
    "I65" = 0
    "J65" = 0
    "K65" = 0
    "L65" = 0
    "H107" = 2624
    "I107" = 2624
    "J107" = 0
    "K107" = 2624
    "L107" = 0
    "M107" = -704
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
foreach ($ref in @("M62","N62","M65","N65","N107")) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H106" = 49990
    "J106" = 49990
    "L106" = 49990
    "N106" = -52514
    "H116" = 114996.664
    "J116" = 114996.664
    "L116" = 114996.664
    "N116" = -124174.664
    "H132" = 1795.2222
    "I132" = 1410.1428
    "K132" = 4230.428400000001
    "M132" = -1700.428400000001
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H4" = 18166688
    "I4" = 18666700
    "K4" = 56000100
    "M4" = -55999988
    "H129" = 55721484
    "I129" = 1045
    "K129" = 3135
    "M129" = 1865
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H114" = 109000
    "J114" = 109000
    "L114" = 109000
    "N114" = -117678
    "H118" = 42564.6
    "J118" = 42564.6
    "L118" = 42564.6
    "N118" = -45878.6
    "H132" = 29414836
    "I132" = 34485690
    "K132" = 103457070
    "M132" = -103454540
    "H141" = 54818
    "J141" = 54818
    "L141" = 54818
    "N141" = -65178
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 54248.65
    "I7" = 3596.9333
    "K7" = 3596.9333
    "M7" = -3484.9333
    "H16" = 2481.6667
    "I16" = 2481.6667
    "K16" = 2481.6667
    "M16" = -2311.6667
    "H93" = 45455736
    "I93" = 50001172
    "K93" = 50001172
    "M93" = -49999924
    "H100" = 3403.3845
    "I100" = 3658.5454
    "K100" = 3658.5454
    "M100" = -3117.5454
    "H126" = 54248.65
    "I126" = 3596.9333
    "K126" = 10790.7999
    "M126" = -8320.7999
    "H132" = 97872.38
    "J132" = 169598.17
    "L132" = 508794.51
    "N132" = -513854.51
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H74" = 24944.8
    "J74" = 24944.8
    "L74" = 24944.8
    "N74" = -26816.8
    "H77" = 24944.8
    "J77" = 24944.8
    "L77" = 74834.39999999999
    "N77" = -84194.39999999999
    "H105" = 96871.664
    "J105" = 96871.664
    "L105" = 96871.664
    "N105" = -103859.664
    "H107" = 25001064
    "I107" = 35715490
    "J107" = 741.5
    "K107" = 107146470
    "L107" = 2224.5
    "M107" = -107144550
    "N107" = -6064.5
    "H113" = 1075.6666
    "I113" = 1058.8
    "K113" = 3176.4
    "M113" = -1006.4
    "H117" = 59996.668
    "J117" = 59996.668
    "L117" = 59996.668
    "N117" = -69174.66800000001
    "H132" = 3997.842
    "I132" = 3892.2666
    "K132" = 11676.7998
    "M132" = -9146.799800000001
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

